$d = $word.ActiveDocument

# Locate the paragraph with the "Added paperclip hole..." text (last bullet
# of the FreedomWing 1.1 changelog section) and append the new list items
# right after it, matching the existing ListParagraph / numId=3 list.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Added paperclip hole to access additional button on Feather RP2040*") {
        $target = $p
        break
    }
}

$items = @(
    @{ Text = "Added documentation"; Level = 1 },
    @{ Text = "3D Printing Guide"; Level = 2 },
    @{ Text = "Assembly Guide"; Level = 2 },
    @{ Text = "Design Rationale"; Level = 2 },
    @{ Text = "Bill of Materials"; Level = 2 },
    @{ Text = "User Guide"; Level = 2 }
)

$current = $target
foreach ($item in $items) {
    $current.Range.InsertParagraphAfter()
    $current = $current.Next()
    $current.Range.Text = $item.Text
    $current.Range.ListFormat.ListLevelNumber = $item.Level
}
